# Apply the update described by the commit "Upload new version with timestamp":
#  - A new product "فازلين بيور كبير" enters the shared-strings table.
#  - Row 25 (item 22) and Row 26 (item 23) pick up new label/quantity data.
#  - A brand-new row (item 24, "معجون اسنان فلورو بالكولا") is inserted at row 27.
#  - The totals row (was row 27) moves to row 28 and its total is recalculated.
#  - The footer row (was row 28) moves to row 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 25 (item 22) ---
$ws.Range("B25").Value = "فازلين بيور كبير"
$ws.Range("H25").Value = "4:0"
$ws.Range("L25").Value = 40

# --- Update existing row 26 (item 23) ---
$ws.Range("B26").Value = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H26").Value = "2:0"
$ws.Range("L26").Value = 20

# --- Insert a brand-new row 27 for item 24, copying formatting from row 26 ---
$ws.Rows.Item(27).Insert()
$ws.Range("A26:N26").Copy()
$ws.Range("A27:N27").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A27").Value = 24
$ws.Range("B27").Value = "معجون اسنان فلورو بالكولا"
$ws.Range("H27").Value = "3:0"
$ws.Range("L27").Value = 30
$ws.Range("N27").Value = "1:0"

$ws.Range("B27:G27").Merge()
$ws.Range("H27:K27").Merge()
$ws.Range("L27:M27").Merge()

$ws.Rows.Item(27).RowHeight = 25.5

# --- Totals row, now shifted down to row 28 ---
$ws.Range("K28").Value = 1596.04
$ws.Rows.Item(28).RowHeight = 25.5

# --- Footer row, now shifted down to row 29 ---
$ws.Rows.Item(29).RowHeight = 17.25
